# Add a new "october-2025" worksheet at the end of the workbook (after
# "september-2025"), mirroring the single-cell "tax revenue" summary sheets
# used for every other month, and record its shared-string value.

$wb = $excel.ActiveWorkbook

# Remember which sheet was active so the user's current view isn't disturbed
# by inserting the new sheet (Worksheets.Add activates the newly added sheet).
$originalActiveSheetName = $wb.ActiveSheet.Name

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "october-2025"

$newSheet.Range("A1").Value = ": tax revenue                                               77,049           76,342             707              0.9%"

$wb.Worksheets.Item($originalActiveSheetName).Activate()
